# "added array questions on 3_12_15"
#
# 1) Paragraph 6 ("6. Accept a number and find its binary equivalent.") was typed
#    as three separate runs with the _GoBack bookmark sitting in the middle of the
#    word "binary" / " equivalent." split. Collapse it back to a single run via
#    Find/Replace (this also removes the now-stale _GoBack bookmark from here).
# 2) A blank spacer paragraph, a bold+underlined "Arrays" heading, and four new
#    array-themed questions (Q1-Q4, with follow-up paragraphs for Q1 and Q2) are
#    appended at the end of the document. The trailing _GoBack bookmark now marks
#    the very end of the new Q4 paragraph, where the author last left the cursor.

$d = $word.ActiveDocument

# --- 1. Collapse the "6. Accept a number and find its binary equivalent." runs ---
$d.Content.Find.Execute('6. Accept a number and find its binary equivalent.', $true, $false, $false, $false, $false, $true, 1, $false, '6. Accept a number and find its binary equivalent.', 2) | Out-Null

# --- 2. Append the new "Arrays" section paragraph by paragraph ---
# blank spacer paragraph
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>') | Out-Null

# "Arrays" section heading (bold + underline)
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Arrays</w:t></w:r></w:p>') | Out-Null

# Q1: reverseIt() prompt
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t>Q1:</w:t></w:r><w:r><w:t xml:space="preserve"> Write a function </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>reverseIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">), that reverse a String. Use </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a for</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> loop that swaps the first and last characters, then second and next-to-last characters and so on. The string should be passed to </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>reverseIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) method as an argument and it should return the reverse string. </w:t></w:r></w:p>') | Out-Null

# Q1 follow-up: "Write a program to exercise reverseIt()..."
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Write a program to </w:t></w:r><w:r><w:t>exercise</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>reverseIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">). The program should get the string from the user. Call </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>reverseIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) and print the output. Check the program with following String “Able was I era I saw </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>elba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>.”</w:t></w:r></w:p>') | Out-Null

# Q2: Employee class prompt
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t>Q2</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Create a class called </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Employee  that</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> contains a name and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>employeeid</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Include the member functions as </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>putData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() respectively for getting and displaying values. </w:t></w:r></w:p>') | Out-Null

# Q2 follow-up: "Write a main class to exercise Employee class..."
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Write a main class to </w:t></w:r><w:r><w:t xml:space="preserve">exercise Employee class. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>it</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> should create the array of type employee and then allow the user to input the data for 10 employees, ask user to enter the employee id and display the information of the entered employee id. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>') | Out-Null

# Q3: maxInt() prompt
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Q3: </w:t></w:r><w:r><w:t xml:space="preserve">Start with a program that allows the user to input a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>number  of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> integers, and then store them in an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> array. Write a function called </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>maxInt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) that goes through the array, element by element, looking for the largest one. The function should take as arguments the address of the array and the number of elements in it, and return the index number of the largest element. The program should call this function and then display the largest element and its index number. </w:t></w:r></w:p>') | Out-Null

# Q4: prompt (ends where the _GoBack bookmark now lives)
$p = $d.Paragraphs.Add()
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Q4: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>') | Out-Null

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
